# Append three new benchmark rows (commit, time, kpoints/sec, notes) to Sheet1.
# Cells are written in the order below so the shared-string table ends up
# indexed in the same sequence as the target workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: a88d2af
$ws.Range("D4").Value = "Change thread count setting"
$ws.Range("A5").Value = "1e8b92d"
$ws.Range("A4").Value = "a88d2af"
$ws.Range("D5").Value = "Make Slice DataDictionary lazy eval"
$ws.Range("D6").Value = "Make algo run on dedicated thread"

$ws.Range("B4").Value = 19.16
$ws.Range("C4").Value = 186

# Row 5: 1e8b92d
$ws.Range("B5").Value = 18.37
$ws.Range("C5").Value = 195

# Row 6: (no commit hash recorded)
$ws.Range("B6").Value = 18.3
$ws.Range("C6").Value = 196

# Match the workbook's recorded selection after the edit
$ws.Range("D6").Select() | Out-Null
